$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.379.04'
$ws.Range("E2").Value = '  -0.53%  '
$ws.Range("D3").Value = '2.453.92'
$ws.Range("E3").Value = '  +0.28%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.51'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.45%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.75'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.43%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.534'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").Value = '2.449.96'
$ws.Range("E9").Value = '  +0.24%  '
$ws.Range("E10").Value = '  -0.48%  '
$ws.Range("E11").Value = '  +2.33%  '
$ws.Range("E12").Value = '  -0.61%  '
$ws.Range("E13").Value = '  -2.64%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.39'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.39%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000176'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.07%  '
$ws.Range("D16").Value = '2.896.31'
$ws.Range("E16").Value = '  +0.15%  '
$ws.Range("D17").Value = '62.071.61'
$ws.Range("E17").Value = '  -0.89%  '
$ws.Range("D18").Value = '2.449.07'
$ws.Range("E18").Value = '  -0.29%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.90'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.91%  '
$ws.Range("E20").Value = '  -2.01%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '328.53'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.38%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.13'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.18%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.97'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -6.82%  '
$ws.Range("E24").Value = '  -0.24%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '65.67'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.52%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.38'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.75%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '592.33'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.19%  '
$ws.Range("D28").Value = '2.571.86'
$ws.Range("E28").Value = '  +0.27%  '
$ws.Range("D29").Value = '0.0₃0961'
$ws.Range("E29").Value = '  -3.80%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.01'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.61%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.44'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.34%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.04'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.14%  '
$ws.Range("E33").Value = '  +0.04%  '
$ws.Range("E34").Value = '  -0.96%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.94'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.58%  '
$ws.Range("E36").Value = '  +0.19%  '
$ws.Range("E37").Value = '  -3.34%  '
$ws.Range("E38").Value = '  +0.50%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '152.51'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.11%  '
$ws.Range("E40").Value = '  -0.13%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.44'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.02%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '43.29'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.52%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.73'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.18%  '
$ws.Range("E44").Value = '  +0.00%  '
$ws.Range("E45").Value = '  -4.14%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.66'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.74%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '142.10'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.19%  '
$ws.Range("E48").Value = '  +1.22%  '
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '0.0₆0249'
$ws.Range("E49").Value = '  +10.24%  '
$ws.Range("B50").Value = 'Hedera'
$ws.Range("C50").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0523'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.54%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.80'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.55%  '
